$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (item id 38956)
$ws.Range("H17").Value = 4448.5
$ws.Range("I17").Value = 1220
$ws.Range("J17").Value = 5524.6665
$ws.Range("K17").Value = 3660
$ws.Range("L17").Value = 16573.9995
$ws.Range("M17").Value = -3492
$ws.Range("N17").Value = -16909.9995

# Row 28 (item id 27772)
$ws.Range("H28").Value = 31250868
$ws.Range("I28").Value = 40000690
$ws.Range("J28").Value = 1501.1428
$ws.Range("K28").Value = 40000690
$ws.Range("L28").Value = 1501.1428
$ws.Range("M28").Value = -40000205
$ws.Range("N28").Value = -2471.1428

# Row 69 (item id 12616)
$ws.Range("H69").Value = 36264.914
$ws.Range("I69").Value = 8750
$ws.Range("J69").Value = 42057.527
$ws.Range("K69").Value = 26250
$ws.Range("L69").Value = 126172.581
$ws.Range("M69").Value = -25376
$ws.Range("N69").Value = -127920.581

# Row 72 (item id 12616)
$ws.Range("H72").Value = 36264.914
$ws.Range("I72").Value = 8750
$ws.Range("J72").Value = 42057.527
$ws.Range("K72").Value = 78750
$ws.Range("L72").Value = 378517.743
$ws.Range("M72").Value = -74382
$ws.Range("N72").Value = -387253.743

# Row 80 (item id 12605)
$ws.Range("H80").Value = 3628.6614
$ws.Range("I80").Value = 3884.4
$ws.Range("J80").Value = 3297.1482
$ws.Range("K80").Value = 11653.2
$ws.Range("L80").Value = 9891.444600000001
$ws.Range("M80").Value = -10655.2
$ws.Range("N80").Value = -11887.4446

# Row 83 (item id 12605)
$ws.Range("H83").Value = 3628.6614
$ws.Range("I83").Value = 3884.4
$ws.Range("J83").Value = 3297.1482
$ws.Range("K83").Value = 34959.6
$ws.Range("L83").Value = 29674.3338
$ws.Range("M83").Value = -29967.6
$ws.Range("N83").Value = -39658.3338

# Row 88 (item id 12608)
$ws.Range("H88").Value = 5978.1113
$ws.Range("I88").Value = 6914.2856
$ws.Range("J88").Value = 2701.5
$ws.Range("K88").Value = 6914.2856
$ws.Range("L88").Value = 2701.5
$ws.Range("M88").Value = -6508.2856
$ws.Range("N88").Value = -3513.5

# Row 91 (item id 12608)
$ws.Range("H91").Value = 5978.1113
$ws.Range("I91").Value = 6914.2856
$ws.Range("J91").Value = 2701.5
$ws.Range("K91").Value = 6914.2856
$ws.Range("L91").Value = 2701.5
$ws.Range("M91").Value = -5510.2856
$ws.Range("N91").Value = -5509.5

# Row 96 (item id 19894)
$ws.Range("H96").Value = 737.3684
$ws.Range("I96").Value = 710.5454999999999
$ws.Range("J96").Value = 774.25
$ws.Range("K96").Value = 2131.6365
$ws.Range("L96").Value = 2322.75
$ws.Range("M96").Value = -758.6364999999996
$ws.Range("N96").Value = -5068.75

# Row 106 (item id 19903)
$ws.Range("H106").Value = 2529.8572
$ws.Range("I106").Value = 2201.5
$ws.Range("K106").Value = 2201.5
$ws.Range("M106").Value = -1570.5

# Row 127 (item id 36114)
$ws.Range("H127").Value = 12740.913
$ws.Range("I127").Value = 1269.8
$ws.Range("K127").Value = 3809.4
$ws.Range("M127").Value = 1150.6

# Row 129 (item id 36115)
$ws.Range("H129").Value = 406710.72
$ws.Range("I129").Value = 524867.6
$ws.Range("K129").Value = 1574602.8
$ws.Range("M129").Value = -1569602.8

# Row 137 (item id 44013)
$ws.Range("H137").Value = 4080.0615
$ws.Range("I137").Value = 3425.22
$ws.Range("K137").Value = 10275.66
$ws.Range("M137").Value = -7725.66

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item id 44147)
$ws.Range("H32").Value = 11240746
$ws.Range("I32").Value = 12052585
$ws.Range("K32").Value = 12052585
$ws.Range("M32").Value = -12052298

# Row 61 (item id 43999)
$ws.Range("H61").Value = 2237.9348
$ws.Range("I61").Value = 1975.4884
$ws.Range("K61").Value = 1975.4884
$ws.Range("M61").Value = -1763.4884

# Row 74 (item id 44000)
$ws.Range("H74").Value = 1624.5714
$ws.Range("I74").Value = 841.7273
$ws.Range("J74").Value = 4495
$ws.Range("K74").Value = 841.7273
$ws.Range("L74").Value = 4495
$ws.Range("M74").Value = 32.27269999999999
$ws.Range("N74").Value = -6243

# Row 77 (item id 44000)
$ws.Range("H77").Value = 1624.5714
$ws.Range("I77").Value = 841.7273
$ws.Range("J77").Value = 4495
$ws.Range("K77").Value = 4208.636500000001
$ws.Range("L77").Value = 22475
$ws.Range("M77").Value = 159.3634999999995
$ws.Range("N77").Value = -31211

# Row 122 (item id 36168)
$ws.Range("H122").Value = 1694
$ws.Range("I122").Value = 1512.9333
$ws.Range("K122").Value = 4538.7999
$ws.Range("M122").Value = -2088.7999

# Row 136 (item id 43999)
$ws.Range("H136").Value = 2237.9348
$ws.Range("I136").Value = 1975.4884
$ws.Range("K136").Value = 5926.4652
$ws.Range("M136").Value = -3376.4652

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (item id 43998)
$ws.Range("H134").Value = 2271.2078
$ws.Range("I134").Value = 1907.2142
$ws.Range("J134").Value = 3241.8572
$ws.Range("K134").Value = 5721.642599999999
$ws.Range("L134").Value = 9725.571599999999
$ws.Range("M134").Value = -3186.642599999999
$ws.Range("N134").Value = -14795.5716

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (item id 44021)
$ws.Range("H58").Value = 2773.5334
$ws.Range("I58").Value = 1231.5238
$ws.Range("K58").Value = 1231.5238
$ws.Range("M58").Value = -1028.5238

# Row 134 (item id 44020)
$ws.Range("H134").Value = 3089.389
$ws.Range("I134").Value = 2231.8096
$ws.Range("K134").Value = 6695.4288
$ws.Range("M134").Value = -4160.4288

# Row 136 (item id 44021)
$ws.Range("H136").Value = 2773.5334
$ws.Range("I136").Value = 1231.5238
$ws.Range("K136").Value = 3694.5714
$ws.Range("M136").Value = -1144.5714

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (item id 4650)
$ws.Range("H4").Value = 28886274
$ws.Range("I4").Value = 36107600
$ws.Range("K4").Value = 108322800
$ws.Range("M4").Value = -108322688

# Row 21 (item id 4667)
$ws.Range("H21").Value = 293.8
$ws.Range("I21").Value = 292.25
$ws.Range("K21").Value = 876.75
$ws.Range("M21").Value = -703.75

# Row 35 (item id 4718)
$ws.Range("H35").Value = 331.125
$ws.Range("I35").Value = 166.33333
$ws.Range("J35").Value = 430
$ws.Range("K35").Value = 498.99999
$ws.Range("L35").Value = 1290
$ws.Range("M35").Value = -210.99999
$ws.Range("N35").Value = -1866

# Row 131 (item id 36060)
$ws.Range("H131").Value = 2666.2
$ws.Range("I131").Value = 1577.9412
$ws.Range("K131").Value = 4733.8236
$ws.Range("M131").Value = 306.1764000000003

$ws = $wb.Worksheets.Item("GSM")
# Row 40 (item id 4113)
$ws.Range("H40").Value = 24999
$ws.Range("J40").Value = 24999
$ws.Range("L40").Value = 24999
$ws.Range("N40").Value = -25301

# Row 55 (item id 4237)
$ws.Range("H55").Value = 9412.571
$ws.Range("J55").Value = 16666.334
$ws.Range("L55").Value = 16666.334
$ws.Range("N55").Value = -17320.334

# Row 97 (item id 19940)
$ws.Range("H97").Value = 403.78946
$ws.Range("I97").Value = 406.54544
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 406.54544
$ws.Range("L97").Value = 400
$ws.Range("M97").Value = 89.45456000000001
$ws.Range("N97").Value = -1392

# Row 132 (item id 44008)
$ws.Range("H132").Value = 2798.6365
$ws.Range("I132").Value = 2551.077
$ws.Range("K132").Value = 7653.231000000001
$ws.Range("M132").Value = -5123.231000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (item id 12565)
$ws.Range("H82").Value = 3955.5
$ws.Range("I82").Value = 2452.5454
$ws.Range("K82").Value = 2452.5454
$ws.Range("M82").Value = -2091.5454

# Row 85 (item id 12565)
$ws.Range("H85").Value = 3955.5
$ws.Range("I85").Value = 2452.5454
$ws.Range("K85").Value = 2452.5454
$ws.Range("M85").Value = -1204.5454

# Row 136 (item id 44060)
$ws.Range("H136").Value = 2461.7317
$ws.Range("I136").Value = 1782.6897
$ws.Range("K136").Value = 5348.0691
$ws.Range("M136").Value = -2798.0691

$ws = $wb.Worksheets.Item("WVR")
# Row 27 (item id 27174)
$ws.Range("H27").Value = 60000
$ws.Range("J27").Value = 60000
$ws.Range("L27").Value = 60000
$ws.Range("N27").Value = -60138

# Row 115 (item id 50000)
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# Row 122 (item id 36208)
$ws.Range("H122").Value = 333039.9
$ws.Range("I122").Value = 2353.5386
$ws.Range("J122").Value = 1407770.6
$ws.Range("K122").Value = 7060.6158
$ws.Range("L122").Value = 4223311.800000001
$ws.Range("M122").Value = -4610.6158
$ws.Range("N122").Value = -4228211.800000001

Write-Output "Applied all Ultros_Profits market data updates"